$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44174
$ws.Range("H2").Value = "Verde"
$ws.Range("J2").Value = 100
$ws.Range("L2").Value = 1100
$ws.Range("M2").Value = 1050
$ws.Range("P2").Value = 1050

# Row 3
$ws.Range("D3").Value = 44159
$ws.Range("J3").Value = 2000
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = 1000
$ws.Range("P3").Value = 1000

# Row 4
$ws.Range("D4").Value = 44441
$ws.Range("J4").Value = 40
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 3000
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 3000

# Row 5
$ws.Range("D5").Value = 44161
$ws.Range("J5").Value = 3000

# Row 6
$ws.Range("D6").Value = 44167
$ws.Range("J6").Value = 140
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = 957
$ws.Range("P6").Value = 957

# Row 7
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 650
$ws.Range("K7").Value = 900
$ws.Range("L7").Value = 1100
$ws.Range("M7").Value = 1008
$ws.Range("P7").Value = 1008

# Row 8
$ws.Range("D8").Value = 44165
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 180
$ws.Range("K8").Value = 800
$ws.Range("L8").Value = 800
$ws.Range("M8").Value = 800
$ws.Range("P8").Value = 800

# Row 9
$ws.Range("D9").Value = 44166
$ws.Range("J9").Value = 285
$ws.Range("L9").Value = 1100
$ws.Range("M9").Value = 1054
$ws.Range("P9").Value = 1054

# Row 11
$ws.Range("D11").Value = 44160
$ws.Range("J11").Value = 1400

# Row 13
$ws.Range("D13").Value = 44162
$ws.Range("J13").Value = 1500
$ws.Range("K13").Value = 1200
$ws.Range("L13").Value = 1200
$ws.Range("M13").Value = 1200
$ws.Range("O13").Value = "Región del Bíobío"
$ws.Range("P13").Value = 1200

# New row 14
$ws.Range("D14").NumberFormat = $ws.Range("D13").NumberFormat()
$ws.Range("A14").Value = 10
$ws.Range("B14").Value = "Vega Modelo de Temuco"
$ws.Range("C14").Value = "La Araucanía"
$ws.Range("D14").Value = 44162
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = 300000000
$ws.Range("G14").Value = "Espárragos"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 1200
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("M14").Value = 1000
$ws.Range("N14").Value = "$/kilo"
$ws.Range("O14").Value = "Región del Maule"
$ws.Range("P14").Value = 1000
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = "Hortaliza"
